# Add the missing data points to the "batch size = 20" and "batch size = 100"
# sheets (each was missing the last "2 machines" accuracy reading), then
# leave the workbook's view state the way the author left it when they saved:
# "batch size = 100" is now the active tab, with the lingering selections at
# C30 (on "batch size = 20") and D20 (on "batch size = 100").

$wb = $excel.ActiveWorkbook

# ---- "batch size = 20" sheet: fill in the missing C16 value ----
$ws20 = $wb.Worksheets.Item("batch size = 20")
$ws20.Range("C16").Value = 0.916

# Refresh the chart that plots column C so it picks up the new point.
$chart20 = $ws20.ChartObjects(1).Chart
$chart20.Refresh()

# Leave the lingering selection on this (now inactive) sheet at C30.
[void]$ws20.Range("C30").Select()

# ---- "batch size = 100" sheet: fill in the missing D15 value ----
$ws100 = $wb.Worksheets.Item("batch size = 100")
$ws100.Range("D15").Value = 0.9172

# Refresh the chart that plots column D so it picks up the new point.
$chart100 = $ws100.ChartObjects(1).Chart
$chart100.Refresh()

# Leave the selection at D20 and make this the active (visible) sheet/tab.
[void]$ws100.Range("D20").Select()
$ws100.Activate()
